# Replace the text sample-number labels in column C ("5d 1", "5d 2", ...,
# "JR3 10") with plain sequential numbers (1..31), matching the new
# "template" format used across the metadata files.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 0; $i -lt 31; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 3).Value = $i + 1
}

# Reset the view: clear the frozen/scrolled top-left cell and move the
# selection from column K to column C.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("C2:C32").Select() | Out-Null
